$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C5").Value = 67704.39484783061
$ws.Range("C6").Value = 67704.39484783061
$ws.Range("C7").Value = 65673.26300239567
$ws.Range("C9").Value = 18846.614647889844
$ws.Range("C10").Value = 48857.78019994077
$ws.Range("C11").Value = 48857.78019994077
$ws.Range("C12").Value = 13780.0
$ws.Range("C14").Value = 35077.780199940775
$ws.Range("C15").Value = 34356.65960238931
$ws.Range("C16").Value = 338.54785505146185
$ws.Range("C18").Value = 1946.1
$ws.Range("C19").Value = 32749.10745744078
$ws.Range("C20").Value = 18334.760616664655

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C2").Value = 6770.957101029237
$ws.Range("C3").Value = 7282.5
$ws.Range("D3").Value = 7.554957022146899
$ws.Range("C5").Value = 7282.499999999999
$ws.Range("D8").Value = 36.361815061514925
$ws.Range("C9").Value = 8179.0
$ws.Range("D9").Value = 20.79533037887259
$ws.Range("C10").Value = 6502.0
$ws.Range("D10").Value = -3.972216881840133
$ws.Range("D11").Value = 16.541869668624965
$ws.Range("C12").Value = 6940.0
$ws.Range("D12").Value = 2.4965879483281266
$ws.Range("C13").Value = 6658.0
$ws.Range("D13").Value = -1.6682589971226705
$ws.Range("D14").Value = -4.57774427461844
$ws.Range("D15").Value = -5.537726726584049

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C2").Value = 7177.2145270909905
$ws.Range("C3").Value = 5492.428571428571
$ws.Range("D3").Value = -23.474092202525853
$ws.Range("C5").Value = 5492.428571428571
$ws.Range("C8").Value = 6054.0
$ws.Range("D8").Value = -15.64972766038027
$ws.Range("C9").Value = 4294.0
$ws.Range("D9").Value = -40.1717757802565
$ws.Range("C10").Value = 894.0
$ws.Range("D10").Value = -87.54391419365378
$ws.Range("C11").Value = 7143.0
$ws.Range("D11").Value = -0.4767103862068497
$ws.Range("C12").Value = 5871.0
$ws.Range("D12").Value = -18.199463345571946
$ws.Range("C13").Value = 8390.0
$ws.Range("D13").Value = 16.897718026000913
$ws.Range("D14").Value = -19.174772077612477

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C2").Value = 778.6600666183622
$ws.Range("C3").Value = 625.7142857142856
$ws.Range("D3").Value = -19.642176022754548
$ws.Range("C5").Value = 625.7142857142856
$ws.Range("C8").Value = 409.0
$ws.Range("D8").Value = -47.4738698523165
$ws.Range("D9").Value = -10.10197774235097
$ws.Range("D10").Value = -38.09879038830452
$ws.Range("C11").Value = 737.0
$ws.Range("D11").Value = -5.3502251373038066
$ws.Range("D12").Value = -39.38304784912808
$ws.Range("C13").Value = 519.0
$ws.Range("D13").Value = -33.34703778325736
$ws.Range("D14").Value = 36.25971659337946

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C2").Value = 778.6600666183622
$ws.Range("C3").Value = 465.66666666666663
$ws.Range("D3").Value = -40.19641090764966
$ws.Range("C5").Value = 465.66666666666663
$ws.Range("D8").Value = -35.53027546665741
$ws.Range("D9").Value = -38.09879038830452
$ws.Range("C10").Value = 498.0
$ws.Range("D10").Value = -36.043978450986835
$ws.Range("D11").Value = -50.9413649965401
$ws.Range("C12").Value = 183.0
$ws.Range("D12").Value = -76.49808846692889
$ws.Range("D13").Value = -4.065967676480249

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C2").Value = 2572.9636983911096
$ws.Range("C3").Value = 2782.0
$ws.Range("D3").Value = 8.124339326652846
$ws.Range("C5").Value = 2781.9999999999995
$ws.Range("D10").Value = 9.601235406599937
$ws.Range("C11").Value = 1374.0
$ws.Range("D11").Value = 6.802905991963343
$ws.Range("C12").Value = 1389.0
$ws.Range("D12").Value = 7.968876581395257
$ws.Range("C14").Value = 1391.0
$ws.Range("D17").Value = 9.601235406599937
$ws.Range("C18").Value = 1374.0
$ws.Range("D18").Value = 6.802905991963343
$ws.Range("C19").Value = 1389.0
$ws.Range("D19").Value = 7.968876581395257
$ws.Range("C21").Value = 1391.0

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("C2").Value = 5619.894393854267
$ws.Range("D3").Value = 14.794565185927684
$ws.Range("D11").Value = 23.062810709808247
$ws.Range("D12").Value = 5.126530606354414
$ws.Range("D13").Value = 16.19435424162057
$ws.Range("D18").Value = 23.062810709808247
$ws.Range("D19").Value = 5.126530606354414
$ws.Range("D20").Value = 16.19435424162057

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C2").Value = 2776.092411421987
$ws.Range("C3").Value = 1686.451092855138
$ws.Range("D3").Value = -39.25090224243314
$ws.Range("C5").Value = 1686.4510928551379
$ws.Range("C9").Value = 1686.451092855138
$ws.Range("D9").Value = -39.25090224243314
$ws.Range("C11").Value = 216.94217584712888
$ws.Range("C13").Value = 1469.5089170080093

$ws = $wb.Worksheets.Item("SYSTEMS")
$ws.Range("C2").Value = 9208.501657399762
$ws.Range("C3").Value = 7963.013507442807
$ws.Range("D3").Value = -13.525415928617502
$ws.Range("C4").Value = 7963.013507442807
$ws.Range("C8").Value = 7963.013507442808
$ws.Range("D8").Value = -13.525415928617482
$ws.Range("C21").Value = 989.3371461547074
$ws.Range("C23").Value = 989.3371461547073
$ws.Range("C26").Value = 491.2608278969228
$ws.Range("C28").Value = 491.26082789692276
$ws.Range("C36").Value = 816.8041180658586
$ws.Range("C38").Value = 816.8041180658585
$ws.Range("C41").Value = 3094.35424026887
$ws.Range("C43").Value = 3094.35424026887
